$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure columns D and E keep their original text representation (prices/percentages
# are stored as literal strings in this workbook, e.g. "6.65%", not as numbers).
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '331.33'
$ws.Range("E2").Value = '6.65%'
$ws.Range("D3").Value = '41.07'
$ws.Range("E3").Value = '9.01%'
$ws.Range("D4").Value = '5.263'
$ws.Range("E4").Value = '1.82%'
$ws.Range("D5").Value = '0.08106'
$ws.Range("E5").Value = '2.46%'
$ws.Range("B6").Value = 'GateToken'
$ws.Range("C6").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D6").Value = '4.526'
$ws.Range("E6").Value = '2.17%'
$ws.Range("B7").Value = 'KuCoinToken'
$ws.Range("C7").Value = 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'
$ws.Range("D7").Value = '8.659'
$ws.Range("E7").Value = '4.51%'
$ws.Range("B8").Value = 'FTXToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D8").Value = '1.910'
$ws.Range("E8").Value = '-0.66%'
$ws.Range("B9").Value = 'BTSEToken'
$ws.Range("C9").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D9").Value = '2.957'
$ws.Range("E9").Value = '-1.44%'
$ws.Range("B10").Value = 'MXToken'
$ws.Range("C10").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D10").Value = '0.9350'
$ws.Range("E10").Value = '-0.39%'
$ws.Range("B11").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C11").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D11").Value = '0.1389'
$ws.Range("E11").Value = '25.26%'
$ws.Range("B12").Value = 'WazirX'
$ws.Range("C12").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D12").Value = '0.1967'
$ws.Range("E12").Value = '0.91%'
$ws.Range("B13").Value = 'MandalaExchangeToken'
$ws.Range("C13").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D13").Value = '0.09170'
$ws.Range("E13").Value = '0.69%'
$ws.Range("B14").Value = 'BitrueCoin'
$ws.Range("C14").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D14").Value = '0.03418'
$ws.Range("E14").Value = '2.62%'
$ws.Range("B15").Value = 'BitMartToken'
$ws.Range("C15").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D15").Value = '0.09570'
$ws.Range("E15").Value = '-0.39%'
$ws.Range("B16").Value = 'BitForexToken'
$ws.Range("C16").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D16").Value = '0.001398'
$ws.Range("E16").Value = '0.29%'
$ws.Range("B17").Value = 'TigerCash'
$ws.Range("C17").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D17").Value = '0.006514'
$ws.Range("E17").Value = '10.50%'
$ws.Range("B18").Value = 'LEO'
$ws.Range("C18").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D18").Value = '3.358'
$ws.Range("E18").Value = '-6.60%'
$ws.Range("D19").Value = '0.3523'
$ws.Range("E19").Value = '3.31%'
$ws.Range("D20").Value = '7.229'
$ws.Range("E20").Value = '12.32%'
$ws.Range("D21").Value = '0.1322'
$ws.Range("E21").Value = '3.20%'
$ws.Range("D22").Value = '0.2569'
$ws.Range("E22").Value = '1.83%'
$ws.Range("D23").Value = '0.04443'
$ws.Range("E23").Value = '1.01%'
$ws.Range("E24").Value = '-0.94%'
$ws.Range("D25").Value = '0.004348'
$ws.Range("E25").Value = '-5.26%'
$ws.Range("E26").Value = '-5.21%'
$ws.Range("D27").Value = '0.0003991'
$ws.Range("E27").Value = '-0.07%'
$ws.Range("D39").Value = '0.02546'
$ws.Range("E39").Value = '13.42%'
$ws.Range("D40").Value = '0.05237'
$ws.Range("E40").Value = '2.42%'
$ws.Range("D41").Value = '0.007625'
$ws.Range("E41").Value = '1.94%'
$ws.Range("D42").Value = '0.1430'
$ws.Range("E42").Value = '5.45%'
$ws.Range("D43").Value = '0.009040'
$ws.Range("E43").Value = '1.89%'
$ws.Range("E44").Value = '1.81%'
$ws.Range("D45").Value = '0.008981'
$ws.Range("E45").Value = '-3.67%'
$ws.Range("D46").Value = '0.00006621'
$ws.Range("E46").Value = '0.40%'
$ws.Range("D47").Value = '0.00000000750'
$ws.Range("E47").Value = '-0.07%'
$ws.Range("D48").Value = '0.003341'
$ws.Range("E48").Value = '16.82%'
$ws.Range("D50").Value = '0.00002100'
$ws.Range("E50").Value = '-0.07%'
$ws.Range("D51").Value = '0.0002000'
$ws.Range("E51").Value = '-0.07%'
